# Add a new "2022" column (column N) to the table, mirroring the layout
# and formatting of column M (2021), and move the selection to N2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column M's formatting (borders, number formats, fonts, etc.)
# into column N by copying the whole M2:M10 block across.
$ws.Range("M2:M10").Copy($ws.Range("N2"))
$excel.CutCopyMode = 0

# Now overwrite the copied values with the real 2022 figures.
$ws.Range("N2").Value = $null
$ws.Range("N3").Value = 2022
$ws.Range("N4").Value = 1434
$ws.Range("N5").Value = 12822
$ws.Range("N6").Value = 3099
$ws.Range("N7").Value = 9722
$ws.Range("N8").Value = 14424
$ws.Range("N9").Value = 5279
$ws.Range("N10").Value = 9145

# Move the selection to N2, matching the recorded cursor position.
$ws.Range("N2").Select()
